# Apply the target edit described by the diff:
# - Clear out the old 4-column table (NO[xlsx]/Food/Usage and Apple/Orange/Fish rows)
# - Replace with a 2-column table: numbers 1-4 in column A, and
#   Story / story-text / Orange / Fish in column B
# - Delete the now-unused columns C and D
# - Update the selected cell to D12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$storyText = "On a breezy afternoon in a small village, young Aanya ran to the open field with her new kite. The kite was blue, with bright yellow stars scattered across it. She had waited all week for the perfect day to fly it, and today seemed just right. The wind was strong, and the sky was clear"

# Clear the whole used range first so stale cells (columns C & D) are removed.
$ws.Cells.Clear()

# Column A: row numbers
$ws.Range("A1").Value = 1
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 3
$ws.Range("A4").Value = 4

# Column B: new content
$ws.Range("B1").Value = "Story"
$ws.Range("B2").Value = $storyText
$ws.Range("B3").Value = "Orange"
$ws.Range("B4").Value = "Fish"

# Update the sheet's selection to match the target (activeCell D12)
$ws.Range("D12").Select()

$wb.Save()
